$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.415.34'
$ws.Range("E2").Value = '''  -2.07%  '
$ws.Range("D3").Value = '''3.472.22'
$ws.Range("E3").Value = '''  -4.47%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '''  +0.16%  '
$ws.Range("D5").Value = '''577.49'
$ws.Range("E5").Value = '''  -4.56%  '
$ws.Range("D6").Value = '''192.01'
$ws.Range("E6").Value = '''  -3.90%  '
$ws.Range("D7").Value = '''0.608'
$ws.Range("E7").Value = '''  -3.05%  '
$ws.Range("D8").Value = '''3.460.88'
$ws.Range("E8").Value = '''  -4.48%  '
$ws.Range("E9").Value = '''  +0.05%  '
$ws.Range("D10").Value = '''0.205'
$ws.Range("E10").Value = '''  -6.43%  '
$ws.Range("D11").Value = '''0.617'
$ws.Range("E11").Value = '''  -4.47%  '
$ws.Range("D12").Value = '''51.42'
$ws.Range("E12").Value = '''  -4.54%  '
$ws.Range("D13").Value = '''0.0000285'
$ws.Range("E13").Value = '''  -6.95%  '
$ws.Range("D14").Value = '''9.11'
$ws.Range("E14").Value = '''  -4.86%  '
$ws.Range("D15").Value = '''4.033.83'
$ws.Range("E15").Value = '''  -4.28%  '
$ws.Range("D16").Value = '''641.67'
$ws.Range("E16").Value = '''  +0.56%  '
$ws.Range("D17").Value = '''69.230.70'
$ws.Range("E17").Value = '''  -2.43%  '
$ws.Range("D18").Value = '''3.464.53'
$ws.Range("E18").Value = '''  -4.61%  '
$ws.Range("E19").Value = '''  -5.53%  '
$ws.Range("E20").Value = '''  -1.85%  '
$ws.Range("D21").Value = '''18.15'
$ws.Range("E21").Value = '''  -4.84%  '
$ws.Range("D22").Value = '''0.941'
$ws.Range("E22").Value = '''  -5.78%  '
$ws.Range("D23").Value = '''17.83'
$ws.Range("E23").Value = '''  -2.49%  '
$ws.Range("D24").Value = '''5.31'
$ws.Range("E24").Value = '''  -1.74%  '
$ws.Range("D25").Value = '''99.03'
$ws.Range("E25").Value = '''  -4.99%  '
$ws.Range("D26").Value = '''4.28'
$ws.Range("E26").Value = '''  -7.70%  '
$ws.Range("E27").Value = '''  -4.76%  '
$ws.Range("E28").Value = '''  -6.43%  '
$ws.Range("D29").Value = '''9.29'
$ws.Range("E29").Value = '''  -4.74%  '
$ws.Range("D30").Value = '''32.33'
$ws.Range("E30").Value = '''  -4.47%  '
$ws.Range("E31").Value = '''  -9.71%  '
$ws.Range("D32").Value = '''6.71'
$ws.Range("E32").Value = '''  -6.69%  '
$ws.Range("D33").Value = '''11.59'
$ws.Range("E33").Value = '''  -5.45%  '
$ws.Range("E34").Value = '''  -6.20%  '
$ws.Range("D35").Value = '''60.93'
$ws.Range("E35").Value = '''  -3.83%  '
$ws.Range("D36").Value = '''3.683.17'
$ws.Range("E36").Value = '''  -8.30%  '
$ws.Range("E37").Value = '''  -0.31%  '
$ws.Range("B38").Value = '''Bittensor'
$ws.Range("C38").Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '''510.69'
$ws.Range("E38").Value = '''  -1.14%  '
$ws.Range("B39").Value = '''PEPE'
$ws.Range("C39").Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '''0.0₃0790'
$ws.Range("E39").Value = '''  -10.48%  '
$ws.Range("D40").Value = '''2.94'
$ws.Range("E40").Value = '''  -3.08%  '
$ws.Range("E41").Value = '''  -1.36%  '
$ws.Range("D42").Value = '''0.371'
$ws.Range("E42").Value = '''  -4.86%  '
$ws.Range("E43").Value = '''  -2.64%  '
$ws.Range("B44").Value = '''InjectiveProtocol'
$ws.Range("C44").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''34.31'
$ws.Range("E44").Value = '''  -6.63%  '
$ws.Range("B45").Value = '''CoreDAO'
$ws.Range("C45").Value = '''https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D45").Value = '''3.50'
$ws.Range("E45").Value = '''  +69.20%  '
$ws.Range("E46").Value = '''  -4.29%  '
$ws.Range("D47").Value = '''3.35'
$ws.Range("E47").Value = '''  -4.59%  '
$ws.Range("E48").Value = '''  -3.41%  '
$ws.Range("E49").Value = '''  -4.91%  '
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '''  -0.21%  '
$ws.Range("E51").Value = '''  -6.05%  '
